$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo "Trizol" -> "TRIzol" (affects every cell using this value, G2:G37)
$ws.Cells.Replace("Trizol", "TRIzol")

# 2. Re-enter the roboticRNAPrep column (H2:H37) as an explicit =FALSE() formula
#    instead of a bare boolean literal.
For ($r = 2; $r -le 37; $r++) {
    $ws.Range("H" + $r).Formula = "=FALSE()"
}

# 3. Make the rnaPrepMethod column (G3:G37) match the formatting already used
#    on G2 (copy/paste the format only).
$ws.Range("G2").Copy()
$ws.Range("G3:G37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Leave the selection on the rnaPrepMethod column, matching the final
#    on-screen cursor position after the cleanup pass.
$ws.Range("G2:G37").Select()
